$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on price cells whose new values would otherwise
# be auto-interpreted by Excel as numbers (losing formatting like trailing zeros).
$textCells = @("D5", "D6", "D8", "D9", "D11", "D14", "D18", "D22", "D23", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "35.404.11"
$ws.Range("E2").Value = "  -3.83%  "
$ws.Range("D3").Value = "1.985.72"
$ws.Range("E3").Value = "  -5.17%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "240.56"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").Value = "0.630"
$ws.Range("E6").Value = "  -3.54%  "
$ws.Range("D8").Value = "56.05"
$ws.Range("E8").Value = "  +3.62%  "
$ws.Range("D9").Value = "59.25"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").Value = "  -3.07%  "
$ws.Range("D11").Value = "0.0722"
$ws.Range("E12").Value = "  -6.48%  "
$ws.Range("E13").Value = "  -3.81%  "
$ws.Range("D14").Value = "14.25"
$ws.Range("E14").Value = "  -4.59%  "
$ws.Range("D15").Value = "2.279.87"
$ws.Range("E15").Value = "  -5.16%  "
$ws.Range("E16").Value = "  -4.34%  "
$ws.Range("D17").Value = "1.982.03"
$ws.Range("E17").Value = "  -5.47%  "
$ws.Range("D18").Value = "17.04"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "35.380.40"
$ws.Range("E19").Value = "  -3.87%  "
$ws.Range("E20").Value = "  -3.93%  "
$ws.Range("D21").Value = "0.0₃0831"
$ws.Range("E21").Value = "  -5.48%  "
$ws.Range("D22").Value = "231.31"
$ws.Range("E22").Value = "  -3.24%  "
$ws.Range("D23").Value = "5.00"
$ws.Range("E23").Value = "  -8.35%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  -6.13%  "
$ws.Range("E26").Value = "  +4.54%  "
$ws.Range("D27").Value = "9.09"
$ws.Range("E27").Value = "  -5.71%  "
$ws.Range("D28").Value = "162.59"
$ws.Range("E28").Value = "  -2.67%  "
$ws.Range("D29").Value = "19.35"
$ws.Range("E29").Value = "  -7.82%  "
$ws.Range("E30").Value = "  -3.91%  "
$ws.Range("D31").Value = "1.13"
$ws.Range("E31").Value = "  -3.18%  "
$ws.Range("D32").Value = "4.74"
$ws.Range("E32").Value = "  -9.01%  "
$ws.Range("D33").Value = "0.0582"
$ws.Range("E33").Value = "  -4.03%  "
$ws.Range("D34").Value = "0.0896"
$ws.Range("E34").Value = "  +8.94%  "
$ws.Range("D35").Value = "4.23"
$ws.Range("E35").Value = "  -10.66%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "1.80"
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "2.22"
$ws.Range("E38").Value = "  -8.71%  "
$ws.Range("D39").Value = "4.84"
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("E40").Value = "  -7.05%  "
$ws.Range("D41").Value = "2.80"
$ws.Range("E41").Value = "  -1.67%  "
$ws.Range("D42").Value = "0.0207"
$ws.Range("E42").Value = "  -5.99%  "
$ws.Range("D43").Value = "1.07"
$ws.Range("E43").Value = "  -7.22%  "
$ws.Range("D44").Value = "0.0881"
$ws.Range("E44").Value = "  -8.44%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.364.19"
$ws.Range("E45").Value = "  -3.66%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "89.62"
$ws.Range("E46").Value = "  -6.94%  "
$ws.Range("D47").Value = "7.39"
$ws.Range("E47").Value = "  -5.24%  "
$ws.Range("D48").Value = "15.35"
$ws.Range("E48").Value = "  -4.18%  "
$ws.Range("D49").Value = "2.88"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").Value = "2.24"
$ws.Range("E50").Value = "  -7.82%  "
$ws.Range("D51").Value = "44.99"
$ws.Range("E51").Value = "  -2.04%  "
